# regen sval data to filter save games
# Update the numeric stat columns (TB, d2S, K, IP, sum) for each row.
# Column F (Win) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.455362044514542, 0.002571899574220771, 0.1494219747398047, 0.4942365360607697, 2.101592454889337)
    3 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    4 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 1133.036916526867, 1138.89110869897)
    5 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 1133.036916526867, 1138.89110869897)
    6 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    7 = @(0.01293466051926884, 0.306821227259698, 22.3905356188092, 10.19245300693656, 32.90274451352472)
    8 = @(0.1190320826869504, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 4.457851494814137)
    9 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 10.19245300693656, 11.37104958465707)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
